$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.242.09'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +0.67%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.797.81'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +2.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  +0.29%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.12'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -2.14%  '

$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("E7").Value = '  +17.44%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3744'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  +10.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.74'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -1.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.145'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +1.93%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07568'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +4.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.49'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("E13").Value = '  +0.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.303'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +2.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.540'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +6.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.792.78'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +2.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001093'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +3.32%  '

$ws.Range("E18").Value = '  +2.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '80.87'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +0.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +0.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.52'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +3.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.364'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +2.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.222.12'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +0.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.81'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +1.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.428'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +1.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.56'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +3.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.01'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -1.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.358'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +1.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.994.60'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +2.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.01'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +2.93%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.237'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -2.81%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.039'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -0.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09444'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +8.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.803'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  -0.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2367'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +12.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.14'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +0.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06334'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +2.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02332'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +2.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.192'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +0.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6575'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +1.13%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.389'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +6.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.480'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  -1.41%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.200'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.13'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +2.50%  '

$ws.Range("E45").Value = '  +0.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6110'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +1.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.787'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -0.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.72'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +2.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.029'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +2.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07125'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +1.82%  '

$ws.Range("E51").Value = '  +0.78%  '
